# Regenerate save_data column G ("K" = strikeouts) using the real per-game
# strikeout counts in place of the previous placeholder "Strike#" values.
# (regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$kValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 2
    6  = 2
    7  = 0
    8  = 3
    9  = 2
    10 = 3
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 0
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 3
    24 = 2
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    32 = 2
    33 = 1
    34 = 3
    35 = 1
    36 = 3
    38 = 2
    39 = 1
    40 = 1
    41 = 2
    42 = 2
    43 = 0
    44 = 1
    45 = 1
    46 = 0
    47 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
